$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the templating placeholders from "ticket" to "tickets" (DOMA-746)
# Row 2: {d.ticket[i].X} -> {d.tickets[i].X}
$ws.Range("A2").Value = "{d.tickets[i].address}"
$ws.Range("B2").Value = "{d.tickets[i].processing}"
$ws.Range("C2").Value = "{d.tickets[i].completed}"
$ws.Range("D2").Value = "{d.tickets[i].canceled}"
$ws.Range("E2").Value = "{d.tickets[i].deferred}"
$ws.Range("F2").Value = "{d.tickets[i].closed}"
$ws.Range("G2").Value = "{d.tickets[i].new_or_reopened}"

# Row 3: {d.ticket[i + 1].X} -> {d.tickets[i + 1].X}
$ws.Range("A3").Value = "{d.tickets[i + 1].address}"
$ws.Range("B3").Value = "{d.tickets[i + 1].processing}"
$ws.Range("C3").Value = "{d.tickets[i + 1].completed}"
$ws.Range("D3").Value = "{d.tickets[i + 1].canceled}"
$ws.Range("E3").Value = "{d.tickets[i + 1].deferred}"
$ws.Range("F3").Value = "{d.tickets[i + 1].closed}"
$ws.Range("G3").Value = "{d.tickets[i + 1].new_or_reopened}"

# Update the active cell / selection on the sheet to D28
$ws.Range("D28").Select()

# Nudge the default column width slightly (8.50390625 -> 8.51171875)
$ws.StandardWidth = 8.51171875
